$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("loginTest")

# Bring the loginTest sheet to the front (matches the new activeTab / tabSelected state).
$ws.Activate()

# Update the stored username/password text in A2/B2 (new shared strings).
$ws.Range("A2").Value = "varun.magadiranganath@centurylink.com"
$ws.Range("B2").Value = "Qwerty1@"

# The two hyperlinks on the row swap which mailto: address they point at
# (A2 now carries the link that used to sit on B2's row, and vice versa).
# The existing Hyperlink objects in this engine are read-only snapshots, so
# clear them and re-add in the desired final order/target.
$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:rajesh.yadav@centurylink.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:Hyde@2317") | Out-Null

# Re-adding hyperlinks resets formatting bookkeeping; restore the original
# "Hyperlink" cell style on both cells.
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("B2").Style = "Hyperlink"

# Leave the cursor on B2, matching the saved selection.
$ws.Range("B2").Select() | Out-Null
